$wb = $excel.ActiveWorkbook

# --- 1. Reorder tabs: move GachaEquipTable to come after ShopEquipTable ---
$gachaEquip = $wb.Worksheets.Item("GachaEquipTable")
$shopEquip  = $wb.Worksheets.Item("ShopEquipTable")
$null = $gachaEquip.Move($null, $shopEquip)

# --- 2. Rewrite GachaEquipTable: insert a "rarity|Int" column between grade and prob,
#        and expand the probability table to 7 rows (header + 6 data rows). ---
$gacha = $wb.Worksheets.Item("GachaEquipTable")

$gacha.Range("A1").Value = "grade|Int"
$gacha.Range("B1").Value = "rarity|Int"
$gacha.Range("C1").Value = "prob|float"

$gacha.Range("A2").Value = 3
$gacha.Range("B2").Value = 2
$gacha.Range("C2").Value = 0.005

$gacha.Range("A3").Value = 3
$gacha.Range("B3").Value = 1
$gacha.Range("C3").Value = 0.025

$gacha.Range("A4").Value = 3
$gacha.Range("B4").Value = 0
$gacha.Range("C4").Value = 0.05

$gacha.Range("A5").Value = 2
$gacha.Range("B5").Value = 0
$gacha.Range("C5").Value = 0.2

$gacha.Range("A6").Value = 1
$gacha.Range("B6").Value = 0
$gacha.Range("C6").Value = 0.32

$gacha.Range("A7").Value = 0
$gacha.Range("B7").Value = 0
$gacha.Range("C7").Value = 0.4

# --- 3. Update ShopEquipTable row 4 (the "50-pack" row becomes a "20-pack" row) ---
$shop = $wb.Worksheets.Item("ShopEquipTable")
$shop.Range("A4").Value = "Equip20"
$shop.Range("B4").Value = 20
$shop.Range("C4").Value = 300

# --- 4. Restore UI state: ShopEquipTable active, selection on A3 ---
$shop.Activate()
$shop.Range("A3").Select()
